$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.012.23"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.831.58"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D5").Formula = "'242.40"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Formula = "'0.6259"
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("D7").Formula = "'0.9991"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Formula = "'0.07599"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").Formula = "'0.2920"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Formula = "'22.54"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Formula = "'0.07714"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.837.97"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Formula = "'4.953"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Formula = "'0.6633"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Formula = "'0.00001023"
$ws.Range("E15").Value = "  +18.66%  "
$ws.Range("D16").Formula = "'82.67"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Formula = "'6.044"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "29.012.09"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Formula = "'226.05"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Formula = "'0.9997"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Formula = "'7.188"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Formula = "'0.9989"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Formula = "'158.40"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Formula = "'8.482"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Formula = "'0.1377"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Formula = "'1.488"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Formula = "'4.092"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Formula = "'4.007"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").Formula = "'0.05236"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Formula = "'1.842"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Formula = "'0.7334"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Formula = "'2.686"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("D37").Value = "1.236.15"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Formula = "'2.755"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").Formula = "'0.01782"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Formula = "'6.317"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Formula = "'0.8966"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Formula = "'0.9990"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Formula = "'101.81"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "1.976.70"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Formula = "'0.00000000123"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").Formula = "'64.15"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Formula = "'0.5100"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Formula = "'8.850"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").Formula = "'0.05741"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("D51").Formula = "'6.665"
$ws.Range("E51").Value = "  -0.63%  "
